# Learning diary edit: complete the "ANGULAR:" module section that follows
# the EXPRESSJS section, mirroring commit "angular completed, updated
# README.md and learning diary file".
#
# Strategy: locate the paragraph ending the EXPRESSJS bullet list ("...JSON
# web tokens to authenticate"), strip the _GoBack bookmark that currently
# sits on it, delete the three empty placeholder paragraphs that follow it,
# then inject the full ANGULAR section (heading, "What I learned:" label and
# its eleven bullet points) as one OOXML fragment - re-creating the
# _GoBack bookmark at the very end of the new content, same as upstream.

$d = $word.ActiveDocument

# --- locate the anchor paragraph -------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*JSON web tokens*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate the 'JSON web tokens' anchor paragraph"
}
$anchor = $d.Paragraphs.Item($anchorIndex)

# --- drop the _GoBack bookmark currently sitting on the anchor paragraph ---------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
    # no pre-existing bookmark - nothing to remove
}

# --- remove the three empty placeholder paragraphs right after the anchor --------
$placeholderStart = $d.Paragraphs.Item($anchorIndex + 1).Range.Start
$placeholderEnd = $d.Paragraphs.Item($anchorIndex + 3).Range.End
$placeholderRange = $d.Range($placeholderStart, $placeholderEnd)
$placeholderRange.Delete()

# --- insert the new ANGULAR section right after the anchor paragraph -------------
$insertPoint = $d.Range($anchor.Range.End - 1, $anchor.Range.End - 1)

$newSectionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:ind w:firstLine="240"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:b/>
      <w:bCs/>
      <w:i/>
      <w:iCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:b/>
      <w:bCs/>
      <w:i/>
      <w:iCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">ANGULAR: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:u w:val="single"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>What I learned:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:ind w:firstLine="0" w:firstLineChars="0"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Understand the Angular fundamentals as well as know the structure of an Angular project</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Know how to use CLI to create components and services</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-I am working with React so when move to Angular framework, the logic is not different but the syntax is tricky for me to remember and understand</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Using two-way data binding with the ngModel directive</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Know some condition syntax in Angular: *ngFor, *ngIf&#8230;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">-Know lifecycle hook: ngOnInit </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Know how to defined routes, a redirect route, and a parameterized route.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Know how to share a service among multiple components</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Use HTTP and add post, put, delete, get method to a service</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Learn how to use observables</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="14"/>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>-Angular is more difficult for me because of its syntax and the flow of code. I need to spend more time on researching about this framework.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newSectionXml)

# InsertXML silently drops an explicit w:firstLine="0" (it reads as "no
# override"), so the "-Understand the Angular fundamentals..." paragraph
# loses its <w:ind w:firstLine="0" .../> half. Restore it by touching the
# paragraph formatting directly - setting both the twips and character-unit
# flavours of the property is what makes the writer re-emit w:firstLine="0"
# alongside w:firstLineChars="0".
$firstBulletIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Understand the Angular fundamentals*") {
        $firstBulletIndex = $i
        break
    }
}
if ($firstBulletIndex -ne -1) {
    $firstBullet = $d.Paragraphs.Item($firstBulletIndex)
    $firstBullet.Format.FirstLineIndent = 0
    $firstBullet.Format.CharacterUnitFirstLineIndent = 0
}
